$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 2000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -5708
$ws.Range("N48").ClearContents()

# Row 56
$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 2000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -5466
$ws.Range("N56").ClearContents()

# Row 92
$ws.Range("H92").Value = 7711.7
$ws.Range("I92").Value = 1467.6666
$ws.Range("J92").Value = 10387.714
$ws.Range("K92").Value = 1467.6666
$ws.Range("L92").Value = 10387.714
$ws.Range("M92").Value = -219.6666
$ws.Range("N92").Value = -12883.714

# Row 125
$ws.Range("H125").Value = 5535.5
$ws.Range("I125").Value = 5535.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 49819.5
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -47359.5

# Row 132
$ws.Range("H132").Value = 305917.94
$ws.Range("I132").Value = 344441.22
$ws.Range("J132").Value = 13141
$ws.Range("K132").Value = 1033323.66
$ws.Range("L132").Value = 39423
$ws.Range("M132").Value = -1030793.66
$ws.Range("N132").Value = -44483

# Row 137
$ws.Range("H137").Value = 4817.4634
$ws.Range("I137").Value = 4607.968
$ws.Range("J137").Value = 5466.9
$ws.Range("K137").Value = 13823.904
$ws.Range("L137").Value = 16400.7
$ws.Range("M137").Value = -11273.904
$ws.Range("N137").Value = -21500.7

# Row 138
$ws.Range("H138").Value = 4685.125
$ws.Range("I138").Value = 2453.2307
$ws.Range("J138").Value = 6212.2104
$ws.Range("K138").Value = 7359.6921
$ws.Range("L138").Value = 18636.6312
$ws.Range("M138").Value = -2219.6921
$ws.Range("N138").Value = -28916.6312

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5585.0454
$ws.Range("I61").Value = 11963.167
$ws.Range("J61").Value = 4577.9736
$ws.Range("K61").Value = 11963.167
$ws.Range("L61").Value = 4577.9736
$ws.Range("M61").Value = -11751.167
$ws.Range("N61").Value = -5001.9736

# Row 97
$ws.Range("H97").Value = 779.8333
$ws.Range("I97").Value = 314.8125
$ws.Range("J97").Value = 4500
$ws.Range("K97").Value = 314.8125
$ws.Range("L97").Value = 4500
$ws.Range("M97").Value = 181.1875
$ws.Range("N97").Value = -5492

# Row 110
$ws.Range("H110").Value = 4316
$ws.Range("I110").Value = 1767.1875
$ws.Range("J110").Value = 6169.6816
$ws.Range("K110").Value = 1767.1875
$ws.Range("L110").Value = 6169.6816
$ws.Range("M110").Value = 277.8125
$ws.Range("N110").Value = -10259.6816

# Row 122
$ws.Range("H122").Value = 2500.8845
$ws.Range("I122").Value = 1778.3636
$ws.Range("J122").Value = 6474.75
$ws.Range("K122").Value = 5335.0908
$ws.Range("L122").Value = 19424.25
$ws.Range("M122").Value = -2885.0908
$ws.Range("N122").Value = -24324.25

# Row 132
$ws.Range("H132").Value = 913288.0600000001
$ws.Range("I132").Value = 1080600.1
$ws.Range("J132").Value = 148433.14
$ws.Range("K132").Value = 3241800.3
$ws.Range("L132").Value = 445299.42
$ws.Range("M132").Value = -3239270.3
$ws.Range("N132").Value = -450359.42

# Row 136
$ws.Range("H136").Value = 5585.0454
$ws.Range("I136").Value = 11963.167
$ws.Range("J136").Value = 4577.9736
$ws.Range("K136").Value = 35889.501
$ws.Range("L136").Value = 13733.9208
$ws.Range("M136").Value = -33339.501
$ws.Range("N136").Value = -18833.9208

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 500487.5
$ws.Range("I7").Value = 666733.3
$ws.Range("J7").Value = 1750
$ws.Range("K7").Value = 666733.3
$ws.Range("L7").Value = 1750
$ws.Range("M7").Value = -666620.3
$ws.Range("N7").Value = -1976

# Row 86
$ws.Range("H86").Value = 3177
$ws.Range("I86").Value = 1398.5652
$ws.Range("J86").Value = 11357.8
$ws.Range("K86").Value = 1398.5652
$ws.Range("L86").Value = 11357.8
$ws.Range("M86").Value = -275.5652
$ws.Range("N86").Value = -13603.8

# Row 89
$ws.Range("H89").Value = 3177
$ws.Range("I89").Value = 1398.5652
$ws.Range("J89").Value = 11357.8
$ws.Range("K89").Value = 6992.826
$ws.Range("L89").Value = 56789
$ws.Range("M89").Value = -1376.826
$ws.Range("N89").Value = -68021

# Row 107
$ws.Range("H107").Value = 6672644.5
$ws.Range("I107").Value = 7148548
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 7148548
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -7146628
$ws.Range("N107").Value = -13840

# Row 134
$ws.Range("H134").Value = 1322186.4
$ws.Range("I134").Value = 1791612.4
$ws.Range("J134").Value = 7793.5
$ws.Range("K134").Value = 5374837.199999999
$ws.Range("L134").Value = 23380.5
$ws.Range("M134").Value = -5372302.199999999
$ws.Range("N134").Value = -28450.5

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1117160.6
$ws.Range("I22").Value = 1275872.8
$ws.Range("J22").Value = 6175.5
$ws.Range("K22").Value = 1275872.8
$ws.Range("L22").Value = 6175.5
$ws.Range("M22").Value = -1275522.8
$ws.Range("N22").Value = -6875.5

# Row 25
$ws.Range("H25").Value = 14000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 14000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 14000
$ws.Range("N25").Value = -14348
$ws.Range("M25").ClearContents()

# Row 31
$ws.Range("H31").Value = 5401.857
$ws.Range("I31").Value = 1292.2
$ws.Range("J31").Value = 7685
$ws.Range("K31").Value = 1292.2
$ws.Range("L31").Value = 7685
$ws.Range("M31").Value = -997.2
$ws.Range("N31").Value = -8275

# Row 34
$ws.Range("H34").Value = 5401.857
$ws.Range("I34").Value = 1292.2
$ws.Range("J34").Value = 7685
$ws.Range("K34").Value = 1292.2
$ws.Range("L34").Value = 7685
$ws.Range("M34").Value = -1090.2
$ws.Range("N34").Value = -8089

# Row 94
$ws.Range("H94").Value = 50003364
$ws.Range("I94").Value = 83335090
$ws.Range("J94").Value = 5781.375
$ws.Range("K94").Value = 83335090
$ws.Range("L94").Value = 5781.375
$ws.Range("M94").Value = -83334639
$ws.Range("N94").Value = -6683.375

# Row 99
$ws.Range("H99").Value = 6175404.5
$ws.Range("I99").Value = 10103379
$ws.Range("J99").Value = 2873.2856
$ws.Range("K99").Value = 10103379
$ws.Range("L99").Value = 2873.2856
$ws.Range("M99").Value = -10101881
$ws.Range("N99").Value = -5869.2856

# Row 107
$ws.Range("H107").Value = 3994.6667
$ws.Range("I107").Value = 3995.6667
$ws.Range("J107").Value = 3993.6667
$ws.Range("K107").Value = 3995.6667
$ws.Range("L107").Value = 3993.6667
$ws.Range("M107").Value = -2075.6667
$ws.Range("N107").Value = -7833.6667

# Row 126
$ws.Range("H126").Value = 6175404.5
$ws.Range("I126").Value = 10103379
$ws.Range("J126").Value = 2873.2856
$ws.Range("K126").Value = 30310137
$ws.Range("L126").Value = 8619.856800000001
$ws.Range("M126").Value = -30307667
$ws.Range("N126").Value = -13559.8568

# Row 134
$ws.Range("H134").Value = 62506484
$ws.Range("I134").Value = 76927624
$ws.Range("J134").Value = 14883.333
$ws.Range("K134").Value = 230782872
$ws.Range("L134").Value = 44649.999
$ws.Range("M134").Value = -230780337
$ws.Range("N134").Value = -49719.999

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 1370.6666
$ws.Range("I80").Value = 1306
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 3918
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -2982
$ws.Range("N80").Value = -6372

# Row 83
$ws.Range("H83").Value = 1370.6666
$ws.Range("I83").Value = 1306
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 11754
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -7074
$ws.Range("N83").Value = -22860

# Row 121
$ws.Range("H121").Value = 44526.43
$ws.Range("I121").Value = 5144.5
$ws.Range("J121").Value = 60279.2
$ws.Range("K121").Value = 15433.5
$ws.Range("L121").Value = 180837.6
$ws.Range("M121").Value = -14123.5
$ws.Range("N121").Value = -183457.6

# Row 129
$ws.Range("H129").Value = 12821675
$ws.Range("I129").Value = 888.5
$ws.Range("J129").Value = 33334934
$ws.Range("K129").Value = 2665.5
$ws.Range("L129").Value = 100004802
$ws.Range("M129").Value = 2334.5
$ws.Range("N129").Value = -100014802

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8929.444
$ws.Range("I70").Value = 9059.5
$ws.Range("J70").Value = 7889
$ws.Range("K70").Value = 9059.5
$ws.Range("L70").Value = 7889
$ws.Range("M70").Value = -8789.5
$ws.Range("N70").Value = -8429

# Row 73
$ws.Range("H73").Value = 8929.444
$ws.Range("I73").Value = 9059.5
$ws.Range("J73").Value = 7889
$ws.Range("K73").Value = 9059.5
$ws.Range("L73").Value = 7889
$ws.Range("M73").Value = -8123.5
$ws.Range("N73").Value = -9761

# Row 102
$ws.Range("H102").Value = 1159901.1
$ws.Range("I102").Value = 2147753.2
$ws.Range("J102").Value = 7406.9443
$ws.Range("K102").Value = 2147753.2
$ws.Range("L102").Value = 7406.9443
$ws.Range("M102").Value = -2146131.2
$ws.Range("N102").Value = -10650.9443

# Row 126
$ws.Range("H126").Value = 27788882
$ws.Range("I126").Value = 50004890
$ws.Range("J126").Value = 18874.875
$ws.Range("K126").Value = 150014670
$ws.Range("L126").Value = 56624.625
$ws.Range("M126").Value = -150012200
$ws.Range("N126").Value = -61564.625

# Row 132
$ws.Range("H132").Value = 50004756
$ws.Range("I132").Value = 76927290
$ws.Range("J132").Value = 5765
$ws.Range("K132").Value = 230781870
$ws.Range("L132").Value = 17295
$ws.Range("M132").Value = -230779340
$ws.Range("N132").Value = -22355

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 6524.778
$ws.Range("I61").Value = 5470.3335
$ws.Range("J61").Value = 8633.666999999999
$ws.Range("K61").Value = 5470.3335
$ws.Range("L61").Value = 8633.666999999999
$ws.Range("M61").Value = -5268.3335
$ws.Range("N61").Value = -9037.666999999999

# Row 113
$ws.Range("H113").Value = 6524.778
$ws.Range("I113").Value = 5470.3335
$ws.Range("J113").Value = 8633.666999999999
$ws.Range("K113").Value = 5470.3335
$ws.Range("L113").Value = 8633.666999999999
$ws.Range("M113").Value = -3300.3335
$ws.Range("N113").Value = -12973.667

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 11890.556
$ws.Range("I122").Value = 5377.9165
$ws.Range("J122").Value = 24915.834
$ws.Range("K122").Value = 16133.7495
$ws.Range("L122").Value = 74747.50199999999
$ws.Range("M122").Value = -13683.7495
$ws.Range("N122").Value = -79647.50199999999

# Row 126
$ws.Range("H126").Value = 5368
$ws.Range("I126").Value = 2479
$ws.Range("J126").Value = 6812.5
$ws.Range("K126").Value = 7437
$ws.Range("L126").Value = 20437.5
$ws.Range("M126").Value = -4967
$ws.Range("N126").Value = -25377.5
